# Updated symbol list on Thu Dec 29 16:40:47 UTC 2022 with GitHub Actions
#
# This script updates the cryptos price list worksheet to reflect refreshed
# prices/volumes pulled from coinranking.com, and reorders three rows
# (BKEXToken / CEJI / KickToken -> KickToken / BKEXToken / CEJI) to match
# the newly scraped ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---------------------------------------------
# Prices are stored as text in the sheet, so a leading apostrophe is used to
# force each numeric-looking value to be written as text rather than a number.
$ws.Range("D2").Value  = "'246.39"
$ws.Range("D3").Value  = "'24.17"
$ws.Range("D4").Value  = "'5.283"
$ws.Range("D5").Value  = "'0.05794"
$ws.Range("D6").Value  = "'6.471"
$ws.Range("D7").Value  = "'3.125"
$ws.Range("D8").Value  = "'0.8183"
$ws.Range("D9").Value  = "'0.8764"
$ws.Range("D10").Value = "'0.1380"
$ws.Range("D11").Value = "'0.06955"
$ws.Range("D12").Value = "'0.03135"
$ws.Range("D13").Value = "'0.02940"
$ws.Range("D14").Value = "'0.09394"
$ws.Range("D15").Value = "'3.740"
$ws.Range("D16").Value = "'0.001527"
$ws.Range("D17").Value = "'0.04707"
$ws.Range("D18").Value = "'0.0005988"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006192"
$ws.Range("D20").Value = "'0.001237"
$ws.Range("D21").Value = "'0.004672"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("D22").Value = "'0.00006098"
$ws.Range("D23").Value = "'3.537"
$ws.Range("D24").Value = "'2.142"
$ws.Range("D25").Value = "'0.3184"
$ws.Range("D26").Value = "'0.1313"
$ws.Range("D28").Value = "'0.0002331"
$ws.Range("D40").Value = "'0.03723"

# --- Rows 41-43 reordering: KickToken, BKEXToken, CEJI ---------------------
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006488"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1056"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003099"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining price (and label) updates -----------------------------------
$ws.Range("D44").Value = "'0.008284"
$ws.Range("D45").Value = "'0.00005275"
$ws.Range("D47").Value = "'0.3599"
$ws.Range("D48").Value = "'0.002298"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D50").Value = "'0.0001999"
